$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Sheet 'Forecast Comparison': shift Week_Start_Date and forecast numbers forward one week ---
$c = $ws.Range("B2")
$c.NumberFormat = "@"
$c.Value = "2025-01-12"
$c.Style = "Normal"
$ws.Range("D2").Value = 138
$ws.Range("E2").Value = 158
$ws.Range("F2").Value = 189
$ws.Range("G2").Value = 222
$ws.Range("H2").Value = 273

$c = $ws.Range("B3")
$c.NumberFormat = "@"
$c.Value = "2025-01-19"
$c.Style = "Normal"
$ws.Range("D3").Value = 197
$ws.Range("E3").Value = 172
$ws.Range("F3").Value = 207
$ws.Range("G3").Value = 246
$ws.Range("H3").Value = 307

$c = $ws.Range("B4")
$c.NumberFormat = "@"
$c.Value = "2025-01-26"
$c.Style = "Normal"
$ws.Range("D4").Value = 134
$ws.Range("E4").Value = 123
$ws.Range("F4").Value = 149
$ws.Range("G4").Value = 178
$ws.Range("H4").Value = 222

$c = $ws.Range("B5")
$c.NumberFormat = "@"
$c.Value = "2025-02-02"
$c.Style = "Normal"
$ws.Range("D5").Value = 133
$ws.Range("E5").Value = 122
$ws.Range("F5").Value = 148
$ws.Range("G5").Value = 176
$ws.Range("H5").Value = 221

$c = $ws.Range("B6")
$c.NumberFormat = "@"
$c.Value = "2025-02-09"
$c.Style = "Normal"
$ws.Range("D6").Value = 136
$ws.Range("E6").Value = 125
$ws.Range("F6").Value = 152
$ws.Range("G6").Value = 182
$ws.Range("H6").Value = 230

$c = $ws.Range("B7")
$c.NumberFormat = "@"
$c.Value = "2025-02-16"
$c.Style = "Normal"
$ws.Range("D7").Value = 136
$ws.Range("E7").Value = 125
$ws.Range("F7").Value = 151
$ws.Range("G7").Value = 182
$ws.Range("H7").Value = 231

$c = $ws.Range("B8")
$c.NumberFormat = "@"
$c.Value = "2025-02-23"
$c.Style = "Normal"
$ws.Range("D8").Value = 131
$ws.Range("E8").Value = 120
$ws.Range("F8").Value = 146
$ws.Range("G8").Value = 178
$ws.Range("H8").Value = 230

$c = $ws.Range("B9")
$c.NumberFormat = "@"
$c.Value = "2025-03-02"
$c.Style = "Normal"
$ws.Range("D9").Value = 139
$ws.Range("E9").Value = 127
$ws.Range("F9").Value = 155
$ws.Range("G9").Value = 190
$ws.Range("H9").Value = 245

$c = $ws.Range("B10")
$c.NumberFormat = "@"
$c.Value = "2025-03-09"
$c.Style = "Normal"
$ws.Range("D10").Value = 136
$ws.Range("E10").Value = 125
$ws.Range("F10").Value = 152
$ws.Range("G10").Value = 183
$ws.Range("H10").Value = 233

$c = $ws.Range("B11")
$c.NumberFormat = "@"
$c.Value = "2025-03-16"
$c.Style = "Normal"
$ws.Range("D11").Value = 134
$ws.Range("E11").Value = 123
$ws.Range("F11").Value = 150
$ws.Range("G11").Value = 182
$ws.Range("H11").Value = 234

$c = $ws.Range("B12")
$c.NumberFormat = "@"
$c.Value = "2025-03-23"
$c.Style = "Normal"
$ws.Range("D12").Value = 127
$ws.Range("E12").Value = 116
$ws.Range("F12").Value = 141
$ws.Range("G12").Value = 174
$ws.Range("H12").Value = 228

$c = $ws.Range("B13")
$c.NumberFormat = "@"
$c.Value = "2025-03-30"
$c.Style = "Normal"
$ws.Range("D13").Value = 128
$ws.Range("E13").Value = 117
$ws.Range("F13").Value = 142
$ws.Range("G13").Value = 175
$ws.Range("H13").Value = 228

$c = $ws.Range("B14")
$c.NumberFormat = "@"
$c.Value = "2025-04-06"
$c.Style = "Normal"
$ws.Range("D14").Value = 127
$ws.Range("E14").Value = 116
$ws.Range("F14").Value = 141
$ws.Range("G14").Value = 172
$ws.Range("H14").Value = 222

$c = $ws.Range("B15")
$c.NumberFormat = "@"
$c.Value = "2025-04-13"
$c.Style = "Normal"
$ws.Range("D15").Value = 121
$ws.Range("E15").Value = 111
$ws.Range("F15").Value = 135
$ws.Range("G15").Value = 166
$ws.Range("H15").Value = 216

$c = $ws.Range("B16")
$c.NumberFormat = "@"
$c.Value = "2025-04-20"
$c.Style = "Normal"
$ws.Range("D16").Value = 122
$ws.Range("E16").Value = 111
$ws.Range("F16").Value = 136
$ws.Range("G16").Value = 166
$ws.Range("H16").Value = 214

$c = $ws.Range("B17")
$c.NumberFormat = "@"
$c.Value = "2025-04-27"
$c.Style = "Normal"
$ws.Range("D17").Value = 116
$ws.Range("E17").Value = 106
$ws.Range("F17").Value = 129
$ws.Range("G17").Value = 158
$ws.Range("H17").Value = 205

# --- Sheet 'Summary': update recomputed metric values ---
$c2 = $ws2.Range("B2")
$c2.NumberFormat = "@"
$c2.Value = "2022-12-25 to 2025-01-05"
$c2.Style = "Normal"

$c2 = $ws2.Range("B4")
$c2.NumberFormat = "@"
$c2.Value = "296"
$c2.Style = "Normal"

$c2 = $ws2.Range("B5")
$c2.NumberFormat = "@"
$c2.Value = "136"
$c2.Style = "Normal"

$c2 = $ws2.Range("B8")
$c2.NumberFormat = "@"
$c2.Value = "14649 units"
$c2.Style = "Normal"

$c2 = $ws2.Range("B9")
$c2.NumberFormat = "@"
$c2.Value = "2154"
$c2.Style = "Normal"

$c2 = $ws2.Range("B10")
$c2.NumberFormat = "@"
$c2.Value = "1143"
$c2.Style = "Normal"

$c2 = $ws2.Range("B11")
$c2.NumberFormat = "@"
$c2.Value = "601"
$c2.Style = "Normal"

$c2 = $ws2.Range("B12")
$c2.NumberFormat = "@"
$c2.Value = "197"
$c2.Style = "Normal"

$c2 = $ws2.Range("B13")
$c2.NumberFormat = "@"
$c2.Value = "2025-01-19"
$c2.Style = "Normal"

$c2 = $ws2.Range("B14")
$c2.NumberFormat = "@"
$c2.Value = "116"
$c2.Style = "Normal"

$c2 = $ws2.Range("B15")
$c2.NumberFormat = "@"
$c2.Value = "2025-04-27"
$c2.Style = "Normal"

